$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column A; this shifts A:H -> B:I mechanically
#    (values, formulas and styles all move one column to the right).
$ws.Columns.Item(1).Insert()

# 2. New column A header + client names (copy style from column B first so the
#    new cells pick up the same cellXf as their row neighbours, then set values).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2:A10").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Value = "Cliente"
$ws.Range("A2").Value = "Estudio Rivarossa"
$ws.Range("A3").Value = "David Berger"
$ws.Range("A4").Value = "Ignacio Zbrun"

# 3. Fill in the new client's row (row 4), which was blank before the edit.
$ws.Range("B4").Value = 20434943966
$ws.Range("C4").Value = 20434943966
$ws.Range("D4").Value = "Merentiel2024"

# Copy formats for the formula cells in row 4 from row 3 (keeps identical cellXf,
# i.e. no new style entries), then set the actual formulas relative to row 4.
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").Copy() | Out-Null
$ws.Range("G4").PasteSpecial(-4122) | Out-Null
$ws.Range("H3").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("I3").Copy() | Out-Null
$ws.Range("I4").PasteSpecial(-4122) | Out-Null

$ws.Range("E4").Formula = '=LEFT(CELL("filename"),FIND("[",CELL("filename"))-1)'
$ws.Range("F4").Formula = '=E4&"Deudas\"'
$ws.Range("G4").Formula = '=IF(B4=B3,1,0)'
$ws.Range("H4").Formula = '=IF(B4=B5,1,0)'
$ws.Range("I4").Formula = '=G4+H4'

# 4. Row 3 gets a fixed (non auto-fit) height, matching the author's edit.
$ws.Range("A3").RowHeight = 14.25

# 5. E5 keeps the "str formula" style (s=3) even though it stays empty - this
#    mirrors the fill-handle artifact left behind when the shared formula
#    range was extended down to row 5 before being cleared again.
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").ClearContents()

# 6. Column widths. Columns B, C, E, F already carry the right bestFit width
#    (they are the untouched old A, B, D, E columns, just shifted right by the
#    insert above). Only the brand new column A and the resized column D
#    (old C, widened to fit "Merentiel2024") need an explicit width.
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 13

# 7. Sheet view: scroll back to A1 and select C4 (matches the saved view state).
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C4").Select()
